# lion-login-data.xlsx - test excel changes
# Rework Sheet1 from a simple 4-row login fixture into a wide 2-row
# (header + value) fixture covering login + loan-purpose + contact +
# address test data. Sheet2's literal text is untouched (its shared-string
# indices merely get renumbered on save because one string was dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Drop the old rows 3 & 4 (their data is superseded by the new row 2) ---
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# --- Row 1: headers ---
$ws.Range("A1").Value = "userName"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "loginbutton"
$ws.Range("D1").Value = "validateloginmessage"

$ws.Range("E1").Value = "AcceptLoginAlert"
$ws.Range("E1").Style = "Normal"
$ws.Range("F1").Value = "LoanPurpose"
$ws.Range("F1").Style = "Normal"
$ws.Range("G1").Value = "PhoneNumber"
$ws.Range("G1").Style = "Normal"
$ws.Range("H1").Value = "OPT_contact Policy"
$ws.Range("H1").Style = "Normal"
$ws.Range("I1").Value = "SaveLoanPurpose"
$ws.Range("I1").Style = "Normal"
$ws.Range("J1").Value = "validateContactMessage"
$ws.Range("J1").Style = "Normal"
$ws.Range("K1").Value = "acceptErrors"
$ws.Range("K1").Style = "Normal"
$ws.Range("L1").Value = "openAdressForm"
$ws.Range("L1").Style = "Normal"
$ws.Range("M1").Value = "Adress Line"
$ws.Range("M1").Style = "Normal"
$ws.Range("N1").Value = "State"
$ws.Range("N1").Style = "Normal"
$ws.Range("O1").Value = "City "
$ws.Range("O1").Style = "Normal"
$ws.Range("P1").Value = "Zip"
$ws.Range("P1").Style = "Normal"
$ws.Range("Q1").Value = "closeAdressForm"
$ws.Range("Q1").Style = "Normal"
$ws.Range("R1").Value = "PrimaryResidence"
$ws.Range("R1").Style = "Normal"

# --- Row 2: values ---
$ws.Range("A2").Value = "venugopal@gmail.com"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "venu1234"
$ws.Range("C2").Value = "submit"

$ws.Range("D2").Value = ""
$ws.Range("D2").Style = "Normal"

$ws.Range("F2").Value = "Refinance"
$ws.Range("H2").Value = "No"
$ws.Range("I2").Value = "save"
$ws.Range("J2").Value = "Contact Phone Number was left blank"
$ws.Range("K2").Value = "continue"
$ws.Range("L2").Value = "adressForm"
$ws.Range("M2").Value = "259 street, #342"
$ws.Range("N2").Value = "Newyork"
$ws.Range("O2").Value = "Stephentown"
$ws.Range("P2").Value = "12168"
$ws.Range("Q2").Value = "submit"
$ws.Range("R2").Value = "primary"

# --- Column widths (approximate the real-Excel "best fit" autosize) ---
$ws.Columns.Item(1).ColumnWidth = 21.71
$ws.Columns.Item(4).ColumnWidth = 20.57
$ws.Columns.Item(5).ColumnWidth = 16.29
$ws.Columns.Item(6).ColumnWidth = 12.43
$ws.Columns.Item(7).ColumnWidth = 14.14
$ws.Columns.Item(8).ColumnWidth = 18
$ws.Columns.Item(9).ColumnWidth = 16.71
$ws.Columns.Item(10).ColumnWidth = 35.14
$ws.Columns.Item(11).ColumnWidth = 11.86
$ws.Columns.Item(12).ColumnWidth = 16.29
$ws.Columns.Item(13).ColumnWidth = 14.86
$ws.Columns.Item(14).ColumnWidth = 12.29
$ws.Columns.Item(15).ColumnWidth = 13.14
$ws.Columns.Item(16).ColumnWidth = 12.43
$ws.Columns.Item(17).ColumnWidth = 16.29
$ws.Columns.Item(18).ColumnWidth = 17.29

# --- View state: scrolled right, selection on S2 ---
$ws.Range("S2").Select()

# --- Workbook window position/size ---
$win = $excel.Windows.Item(1)
$win.Left = 10545
$win.Top = 135
$win.Width = 13725
$win.Height = 8835
